$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.516.80"
$ws.Range("E2").Value = "  +2.17%  "

# Row 3
$ws.Range("D3").Value = "2.288.84"
$ws.Range("E3").Value = "  +1.54%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "157.32"
$ws.Range("E5").Value = "  +15,603.43%  "

# Row 6
$ws.Range("D6").Value = "306.58"
$ws.Range("E6").Value = "  +1.14%  "

# Row 7
$ws.Range("D7").Value = "97.02"
$ws.Range("E7").Value = "  +6.49%  "

# Row 8
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  +0.58%  "

# Row 9
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("D10").Value = "0.495"
$ws.Range("E10").Value = "  +3.67%  "

# Row 11
$ws.Range("D11").Value = "36.23"
$ws.Range("E11").Value = "  +13.00%  "

# Row 12
$ws.Range("D12").Value = "0.0805"
$ws.Range("E12").Value = "  +1.30%  "

# Row 13
$ws.Range("E13").Value = "  -1.71%  "

# Row 14
$ws.Range("D14").Value = "6.73"
$ws.Range("E14").Value = "  +2.66%  "

# Row 15
$ws.Range("D15").Value = "2.641.14"
$ws.Range("E15").Value = "  +1.45%  "

# Row 16
$ws.Range("D16").Value = "14.59"
$ws.Range("E16").Value = "  +3.17%  "

# Row 17
$ws.Range("D17").Value = "2.265.58"
$ws.Range("E17").Value = "  +0.56%  "

# Row 19
$ws.Range("D19").Value = "42.398.73"
$ws.Range("E19").Value = "  +2.09%  "

# Row 20
$ws.Range("D20").Value = "12.78"
$ws.Range("E20").Value = "  +4.04%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0920"
$ws.Range("E21").Value = "  +2.10%  "

# Row 22
$ws.Range("E22").Value = "  +2.19%  "

# Row 23
$ws.Range("D23").Value = "67.87"
$ws.Range("E23").Value = "  +2.00%  "

# Row 24
$ws.Range("E24").Value = "  +1.43%  "

# Row 25
$ws.Range("D25").Value = "2.61"
$ws.Range("E25").Value = "  +1.26%  "

# Row 26
$ws.Range("E26").Value = "  +2.44%  "

# Row 27
$ws.Range("E27").Value = "  -0.32%  "

# Row 28
$ws.Range("D28").Value = "24.03"
$ws.Range("E28").Value = "  +0.62%  "

# Row 29
$ws.Range("D29").Value = "36.19"
$ws.Range("E29").Value = "  +5.87%  "

# Row 30
$ws.Range("D30").Value = "9.61"
$ws.Range("E30").Value = "  +1.47%  "

# Row 31
$ws.Range("E31").Value = "  +2.01%  "

# Row 32
$ws.Range("D32").Value = "161.91"
$ws.Range("E32").Value = "  +0.54%  "

# Row 33
$ws.Range("D33").Value = "5.33"
$ws.Range("E33").Value = "  +3.72%  "

# Row 34
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.06%  "

# Row 36
$ws.Range("E36").Value = "  +3.47%  "

# Row 37
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "17.35"
$ws.Range("E37").Value = "  +4.81%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.109"
$ws.Range("E38").Value = "  +4.88%  "

# Row 39
$ws.Range("E39").Value = "  +5.19%  "

# Row 40
$ws.Range("E40").Value = "  +0.35%  "

# Row 41
$ws.Range("E41").Value = "  -0.11%  "

# Row 42
$ws.Range("D42").Value = "4.21"
$ws.Range("E42").Value = "  +8.07%  "

# Row 43
$ws.Range("E43").Value = "  +13.40%  "

# Row 44
$ws.Range("D44").Value = "2.003.80"

# Row 45
$ws.Range("D45").Value = "19.39"
$ws.Range("E45").Value = "  -1.18%  "

# Row 46
$ws.Range("E46").Value = "  +2.99%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "3.03"
$ws.Range("E47").Value = "  +6.31%  "

# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "10.26"
$ws.Range("E48").Value = "  +1.00%  "

# Row 49
$ws.Range("D49").Value = "54.21"
$ws.Range("E49").Value = "  +5.32%  "

# Row 50
$ws.Range("E50").Value = "  +1.93%  "

# Row 51
$ws.Range("D51").Value = "72.77"
$ws.Range("E51").Value = "  +0.25%  "
